$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The export window rolled forward one day: drop the oldest date row (2025-10-20)
# by deleting the entire row, which shifts every subsequent row up by one.
$ws.Rows.Item(2).Delete()

# The two most recent days (now rows 2 and 3, i.e. 2025-10-21 and 2025-10-22)
# don't have validated video-indexing data yet, so their "No video indexed" /
# "Video indexed" counts are blanked out.
$ws.Range("B2").Value = ""
$ws.Range("C2").Value = ""
$ws.Range("B3").Value = ""
$ws.Range("C3").Value = ""
